$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- sheet1 ("#system") content edits -------------------------------------

# column P ("number" list): P15 value changes and a new P16 value is added
$ws.Range("P15").Value2 = "roundTo(var,closestDigit)"
$ws.Range("P16").Value2 = "whole(var)"

# column X ("web" list): insert a new entry at X17, pushing X17:X123 down to
# X18:X124 (only column X shifts; the other columns on those rows must stay
# put, so shift values manually instead of using Range.Insert which moves
# the whole row)
for ($r = 123; $r -ge 17; $r--) {
    $ws.Cells.Item($r + 1, 24).Value2 = $ws.Cells.Item($r, 24).Value2
}
$ws.Range("X17").Value2 = "assertElementsPresent(prefix)"

# extend the sheet's recorded dimension down to row 124 (keeping the
# existing, already-stale "AD" column boundary) without introducing a new
# style record
$ws.Range("AD124").Style = "Normal"

# --- workbook-level defined names ------------------------------------------
$wb.Names.Item("number").RefersTo = "='#system'!`$P`$2:`$P`$16"
$wb.Names.Item("web").RefersTo = "='#system'!`$X`$2:`$X`$124"
